# Adds SVM hyperparameter-tuning results (three C/Kernel/Gamma bullet
# groups, each followed by Validation/Test dataset lines and a blank
# spacer line) right after the "SVM:" heading paragraph.

$d = $word.ActiveDocument

# --- Locate the "SVM:" paragraph -------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("SVM:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$svmPara = $rng.Paragraphs(1)
$cursorIdx = $svmPara.Index

# --- Grab a bullet-list template to clone for the new numId ----------------------
# (paragraph 6 is one of the existing "Random_state = 0" bullet items)
$bulletTemplate = $d.Paragraphs(6).Range.ListFormat.ListTemplate
$newListTemplate = $null

# --- Content to insert, in order --------------------------------------------------
$items = @(
    @{ Text = "C = 1000.0"; Bullet = $true },
    @{ Text = "Kernel = rbf"; Bullet = $true },
    @{ Text = "Gamma = auto"; Bullet = $true },
    @{ Text = "Validation dataset = 0.7647058823529411"; Bullet = $false },
    @{ Text = "Test dataset = 0.6417910447761194"; Bullet = $false },
    @{ Text = ""; Bullet = $false },

    @{ Text = "C = 1000.0"; Bullet = $true },
    @{ Text = "Kernel = linear"; Bullet = $true },
    @{ Text = "Gamma = linear"; Bullet = $true },
    @{ Text = "Validation dataset = 0.75"; Bullet = $false },
    @{ Text = "Test dataset = 0.6417910447761194"; Bullet = $false },
    @{ Text = ""; Bullet = $false },

    @{ Text = "C = 100.0"; Bullet = $true },
    @{ Text = "Kernel = rbf"; Bullet = $true },
    @{ Text = "Gamma = auto"; Bullet = $true },
    @{ Text = "Validation dataset = 0.6323529411764706"; Bullet = $false },
    @{ Text = "Test dataset = 0.5373134328358209"; Bullet = $false }
)

foreach ($item in $items) {
    $cursorPara = $d.Paragraphs($cursorIdx)
    $cursorPara.Range.InsertParagraphAfter()
    $cursorIdx = $cursorIdx + 1
    $newPara = $d.Paragraphs($cursorIdx)

    if ($item.Text -ne "") {
        $newPara.Range.Text = $item.Text
    }

    if ($item.Bullet) {
        if ($null -eq $newListTemplate) {
            $newPara.Range.ListFormat.ApplyListTemplate($bulletTemplate)
            $newListTemplate = $newPara.Range.ListFormat.ListTemplate
        } else {
            $newPara.Range.ListFormat.ApplyListTemplate($newListTemplate)
        }
    } else {
        $newPara.Format.LeftIndent = 72
    }
}

Write-Output "Inserted $($items.Count) paragraphs after 'SVM:' (now $($d.Paragraphs.Count) total)."
